# Regenerate save_data column G (header "K" = strikeouts) to use the
# actual strikeout count (K) instead of the previous "Strike#" derived
# value. The new values are written directly into column G, row by row,
# matching the regenerated s_vals output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 4
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 3
    14 = 0
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 0
    21 = 2
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 2
    39 = 1
    40 = 0
    41 = 3
    42 = 2
    43 = 0
    44 = 2
    45 = 1
    46 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
